# Applies the StructureDefinition metadata refresh (IBM/Alvearie -> LinuxForHealth)
# plus the removal of the stale ele-1/ext-1 constraint text from the root
# Extension element row on the Elements sheet.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/rx-formulary-indicator"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet --------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the root "Extension" element; clear its Constraint(s) value (column AI)
$elements.Range("AI2").Value = ""
